# Update "Förändrad" date column (C) for rows 2-5 from 2023-09-01 (45170)
# to 2023-09-05 (45174) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value = 45174
